# Insert 5 new weekly-report rows before the existing row 814, pushing the
# previously-last rows (814-831) down to (819-836), then populate the newly
# inserted rows 814-818 with the new week's data (date 44448).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows above row 814 (existing rows 814-831 shift to 819-836).
$ws.Rows("814:818").Insert()

# New data for the 5 inserted rows.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F CategoriaID,
#          G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo,
#          L Precio maximo, M Precio promedio ponderado, N Unidad,
#          O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificacion

$rows = @(
    @{ R=814; H="Zafiro amarillo"; I="Primera"; J=95;  K=45000; L=50000; M=46579; P=3105 },
    @{ R=815; H="Zafiro rojo";     I="Extra";   J=155; K=50000; L=50000; M=50000; P=3333 },
    @{ R=816; H="Zafiro rojo";     I="Primera"; J=215; K=45000; L=45000; M=45000; P=3000 },
    @{ R=817; H="Zafiro verde";    I="Extra";   J=125; K=45000; L=45000; M=45000; P=3000 },
    @{ R=818; H="Zafiro verde";    I="Primera"; J=375; K=40000; L=43000; M=41520; P=2768 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 4).Value = 44448
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = 100112002
    $ws.Cells.Item($r, 7).Value = "Pimiento"
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = "$/caja 15 kilos"
    $ws.Cells.Item($r, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = 15
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
